$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username test-data value for the registration mass data
$ws.Range("B2").Value = "BRUNO370"

# Move the active selection to B2 (matches the saved view state)
$ws.Range("B2").Select()

$wb.Save()
